# Make plot and visit (settings) forms more consistent:
# - rename the "name" field reference to "plot_name" on the survey sheet
#   and the settings sheet (adds a new shared string "plot_name")
# - update the active sheet / selections to match the new workflow

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# survey!D3 used to reference "name" (shared string 2); it now references
# the new "plot_name" string.
$survey.Range("D3").Value = "plot_name"
$survey.Range("D24").Select() | Out-Null

# settings!B5 used to reference "name" (shared string 2); it now references
# the new "plot_name" string, and the settings sheet becomes the active tab.
$settings.Range("B5").Value = "plot_name"
$settings.Activate() | Out-Null
$settings.Range("B5").Select() | Out-Null
